$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append three new rows of data (171-173) with date labels (column A),
# and values in columns C and D, matching the pattern of existing rows.
$ws.Range("A171").Value = "20-09-2021"
$ws.Range("C171").Value = 2.95
$ws.Range("D171").Value = 2.91

$ws.Range("A172").Value = "23-09-2021"
$ws.Range("C172").Value = 2.73
$ws.Range("D172").Value = 3.27

$ws.Range("A173").Value = "24-09-2021"
$ws.Range("C173").Value = 2.79
$ws.Range("D173").Value = 3.13
